$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.043
$ws.Range("D3").Value = -7.303
$ws.Range("E8").Value = 16.764
$ws.Range("E11").Value = 17.103
$ws.Range("A12").Value = -21.683
$ws.Range("C14").Value = -12.373
$ws.Range("E14").Value = 17.351
$ws.Range("E15").Value = 15.927
$ws.Range("E17").Value = 16.787
$ws.Range("D20").Value = -7.755
$ws.Range("D25").Value = -8.061000000000002
$ws.Range("C26").Value = -12.117
$ws.Range("E26").Value = 16.612
$ws.Range("A27").Value = -21.814
$ws.Range("D30").Value = -7.206
$ws.Range("C31").Value = -12.482
$ws.Range("A32").Value = -21.751
$ws.Range("C35").Value = -12.589
$ws.Range("A36").Value = -20.178
$ws.Range("E36").Value = 16.44
$ws.Range("C37").Value = -13.665
$ws.Range("A38").Value = -19.741
$ws.Range("D44").Value = -7.746
$ws.Range("C45").Value = -12.441
$ws.Range("A46").Value = -21.869
$ws.Range("D47").Value = -7.617999999999999
$ws.Range("C52").Value = -11.363
$ws.Range("A54").Value = -22.15
$ws.Range("A55").Value = -22.174
$ws.Range("A56").Value = -21.997
$ws.Range("C57").Value = -13.829
$ws.Range("D58").Value = -8.062000000000001
$ws.Range("E64").Value = 17.185
$ws.Range("A67").Value = -21.534
$ws.Range("A69").Value = -21.544
$ws.Range("A72").Value = -21.481
$ws.Range("D78").Value = -7.803
$ws.Range("E79").Value = 17.266
$ws.Range("C81").Value = -13.243
$ws.Range("A83").Value = -21.636
$ws.Range("C83").Value = -12.809
$ws.Range("D84").Value = -8.294
$ws.Range("A86").Value = -22.264
$ws.Range("D89").Value = -7.234
$ws.Range("E89").Value = 17.078
$ws.Range("A91").Value = -21.587
$ws.Range("D91").Value = -6.910000000000001
$ws.Range("D92").Value = -7.007000000000001
$ws.Range("A93").Value = -21.547
$ws.Range("D96").Value = -7.472
$ws.Range("A99").Value = -20.037
$ws.Range("C100").Value = -12.156
$ws.Range("C102").Value = -13.419
$ws.Range("D102").Value = -7.334999999999999
